$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.640.42"
$ws.Range("E2").Value = "  -3.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.091.47"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.46"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5162"
$ws.Range("E7").Value = "  -1.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4396"
$ws.Range("E8").Value = "  -2.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09211"
$ws.Range("E9").Value = "  +2.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.09"
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.177"
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.05"
$ws.Range("E12").Value = "  +2.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.090.01"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.752"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.171"
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "100.12"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001154"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.009"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.04"
$ws.Range("E19").Value = "  +8.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06645"
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.179"
$ws.Range("E22").Value = "  -1.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.694.27"
$ws.Range("E23").Value = "  -3.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.65"
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("E25").Value = "  -3.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.331.15"
$ws.Range("E26").Value = "  -1.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.86"
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.60"
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.529"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.42"
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.138"
$ws.Range("E31").Value = "  -3.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1048"
$ws.Range("E32").Value = "  -2.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.629"
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("E34").Value = "  -2.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.956"
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.071"
$ws.Range("E36").Value = "  +2.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.28"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02571"
$ws.Range("E38").Value = "  -2.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06717"
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2244"
$ws.Range("E40").Value = "  -2.85%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.41"
$ws.Range("E41").Value = "  -1.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6859"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.287"
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6693"
$ws.Range("E44").Value = "  +4.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.20"
$ws.Range("E45").Value = "  -3.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.305"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.613"
$ws.Range("E47").Value = "  -4.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.218"
$ws.Range("E48").Value = "  -2.30%  "
$ws.Range("E49").Value = "  -4.56%  "
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.165"
$ws.Range("E51").Value = "  -2.27%  "
